$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reindex the crosstab output values to match the desired cut-and-paste order.
$ws.Range("B2").Value = 0.2760736196319019
$ws.Range("C2").Value = 0.2699386503067485
$ws.Range("D2").Value = 0.006134969325153374

$ws.Range("B3").Value = 0.2638036809815951
$ws.Range("C3").Value = 0.1840490797546012
$ws.Range("D3").Value = 0
